# Scenario 4.xlsx edit:
# - Decrease probability of fetal death (p_fetaldeath) for weeks 0-3 (conception weeks)
#   on the "Phase1" and "Phase2" sheets, letting the dependent "continued pregnancy"
#   formulas recalc automatically.
# - Update the active sheet / selection to reflect the state left behind by the author
#   after making the edits (Phase1 becomes the active tab with B2:B5 selected; Phase2's
#   previous scroll position / selection is reset to C2:C5, no longer the active tab).

$wb = $excel.ActiveWorkbook

# --- Phase1: column B (p_fetaldeath) rows 2-5 ---
$ws1 = $wb.Worksheets.Item("Phase1")
$ws1.Activate() | Out-Null

$ws1.Range("B2").Value = 0.1
$ws1.Range("B3").Value = 0.1
$ws1.Range("B4").Value = 0.05
$ws1.Range("B5").Value = 0.05

# Leave the sheet in the state the author left it in: B2:B5 selected, B2 active.
$ws1.Range("B2:B5").Select() | Out-Null

# --- Phase2: column C (p_fetaldeath) rows 2-5 ---
$ws2 = $wb.Worksheets.Item("Phase2")
$ws2.Activate() | Out-Null

$ws2.Range("C2").Value = 0.1
$ws2.Range("C3").Value = 0.1
$ws2.Range("C4").Value = 0.05
$ws2.Range("C5").Value = 0.05

# Reset the selection/scroll position left on this sheet (no longer the active tab).
$ws2.Range("C2:C5").Select() | Out-Null

# Re-activate Phase1 last, since it is the tab that ends up selected/active.
$ws1.Activate() | Out-Null
$ws1.Range("B2:B5").Select() | Out-Null
